# Word COM-interop script implementing the pin_order.docx edit:
#  1. Remove the stray "_GoBack" bookmark that sat right after "Order Form".
#  2. In the pin-order price table, change six "$5.10" prices to "$5.04"
#     (Break "50"/"60"/"70"/"80"/"90"/"100" rows).
#  3. Change the Eagle Pin price from "$5.40" to "$5.32".
#  4. Insert "12% " before "taxes included)" note, and re-create the
#     "_GoBack" bookmark right after the newly inserted text (Word leaves
#     _GoBack at the most-recently-edited spot).

$d = $word.ActiveDocument

# --- 1. Drop the old _GoBack bookmark (it was next to "Order Form"). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Six "5.10" -> "5.04" price edits in the table. ---
# Each match is "5.10" inside a run that also contains nothing else (the
# leading "$" lives in its own separate run already). Word splits the run
# at the "5." / "04" boundary when it is typed over, so we recreate that by
# forcing a format split (toggle bold on/off) around just the last two
# characters before overwriting them.
$rng = $d.Content
for ($i = 0; $i -lt 10; $i++) {
    $found = $rng.Find.Execute("5.10", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }
    $start = $rng.Start
    $end = $rng.End
    $tail = $d.Range($start + 2, $end)
    $tail.Font.Bold = $true
    $tail.Text = "04"
    $fix = $d.Range($start + 2, $start + 4)
    $fix.Font.Bold = $false
    $rng = $d.Range($start + 4, $start + 4)
}

# --- 3. "5.40" -> "5.32" (Eagle Pin row). ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("5.40", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $start2 = $rng2.Start
    $end2 = $rng2.End
    $tail2 = $d.Range($start2 + 2, $end2)
    $tail2.Font.Bold = $true
    $tail2.Text = "32"
    $fix2 = $d.Range($start2 + 2, $start2 + 4)
    $fix2.Font.Bold = $false
}

# --- 4. Insert "12% " before "taxes included)" and re-drop _GoBack there. ---
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("(taxes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $insertAt = $rng3.Start + 1
    $ins = $d.Range($insertAt, $insertAt)
    $ins.InsertBefore("12% ")
    $newRun = $d.Range($insertAt, $insertAt + 4)
    $newRun.Font.Bold = $true
    $newRun.Font.Bold = $false
    $bmPos = $insertAt + 4
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
